# Carga 201511 Noviembre de 2015
#
# Insert a new "Fecha Servicios Inicio" column before the existing
# "Fecha Servicios Fin" column (old column AB / index 28), shifting all
# subsequent columns one place to the right, and append five new trailing
# columns (Tipo de Servicio, Tipo de Cobro, Precio, kilos Integrados,
# Kilo Excedido) after the current last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new column and set its header -----------------------------
$ws.Columns.Item(28).Insert()
$ws.Cells.Item(1, 28).Value = "Fecha Servicios Inicio"

# --- Append the five new trailing header columns ---------------------------
$ws.Cells.Item(1, 51).Value = "Tipo de Servicio"
$ws.Cells.Item(1, 52).Value = "Tipo de Cobro"
$ws.Cells.Item(1, 53).Value = "Precio"
$ws.Cells.Item(1, 54).Value = "kilos Integrados"
$ws.Cells.Item(1, 55).Value = "Kilo Excedido"

# Match the header formatting (bold white font on the blue fill) used by the
# rest of row 1 by copying the format from the neighbouring header cell.
$ws.Range("AX1").Copy()
$ws.Range("AY1:BC1").PasteSpecial(-4122)

# --- Column widths ----------------------------------------------------------
# Reproduce (as closely as this engine's width quantization allows) the
# auto-fit column widths Excel computed for the new columns.
$ws.Columns.Item(28).ColumnWidth = 27.16666666666667
$ws.Columns.Item(51).ColumnWidth = 20.16666666666667
$ws.Columns.Item(52).ColumnWidth = 17.66666666666667
$ws.Columns.Item(53).ColumnWidth = 8.166666666666666
$ws.Columns.Item(54).ColumnWidth = 20
$ws.Columns.Item(55).ColumnWidth = 17
